# Applies the marksheet-update edit:
#  - Row 10/11/12: fill in the Right/Wrong/Not-Attempt/Max/Total summary values
#    and give the label cells (A10/A11/A12) the same "mtitleStyle" as A9.
#  - Row 15..40: the sheet used to carry three repeated (Student Ans / Correct Ans)
#    column-pairs (A/B, D/E, G/H). The G/H pair is dropped entirely, the D/E pair
#    is kept only for rows 16-18, and the "Student Ans" (A) column is populated for
#    most rows with either the matching option (graded correct -> correctStyle) or
#    a different option (graded incorrect -> incorrectStyle). A few rows (25, 38, 40)
#    are intentionally left untouched/blank, matching upstream's partial re-grade.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ----------------------------------------

# A10/A11/A12 need the same style as the other header/title cells (A9) --
# copy formatting from A9 so we reuse the existing "mtitleStyle" cellXf
# instead of minting a new one.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "75/112"

# ---- Drop the 3rd (G/H) answer-pair entirely ---------------------------
$ws.Range("G15:H40").Clear()

# ---- Drop the D/E answer-pair for every row except 16-18 ---------------
$ws.Range("D19:E40").Clear()

# ---- Fill in the D/E pair for rows 16-18 (D mirrors E's "Correct Ans") -
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"

# ---- Populate the "Student Ans" (A) column for rows 16-40 --------------
# correctStyle rows (student answer == correct answer)
$correctRows = 16,17,18,19,20,21,22,23,26,27,29,30,32,33,35,37,39
# incorrectStyle rows (student answer != correct answer)
$incorrectRows = 24,28,31,34,36

$correctAnswers = @{
  16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C"; 20 = "Option B";
  21 = "Option C"; 22 = "Option D"; 23 = "Option D"; 26 = "Option C"; 27 = "Option A";
  29 = "Option D"; 30 = "Option B"; 32 = "Option C"; 33 = "Option D"; 35 = "Option D";
  37 = "Option A"; 39 = "Option D"
}
$incorrectAnswers = @{
  24 = "Option C"; 28 = "Option C"; 31 = "Option B"; 34 = "Option C"; 36 = "Option D"
}

# B10 already carries "correctStyle" (s=5) - reuse it for every correct-answer cell.
$ws.Range("B10").Copy()
foreach ($r in $correctRows) {
  $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

# C10 already carries "incorrectStyle" (s=6) - reuse it for every wrong-answer cell.
$ws.Range("C10").Copy()
foreach ($r in $incorrectRows) {
  $ws.Cells.Item($r, 1).PasteSpecial(-4122)
}

foreach ($r in $correctRows) {
  $ws.Cells.Item($r, 1).Value = $correctAnswers[$r]
}
foreach ($r in $incorrectRows) {
  $ws.Cells.Item($r, 1).Value = $incorrectAnswers[$r]
}
